# Update "想去人数" (want-to-go count) values on both the "展览" sheet
# and the "全部类型" sheet, which contain duplicated data rows.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 5262
    $ws.Range("F3").Value = 162
    $ws.Range("F4").Value = 914
}
